# Applies the "Added a new ancestral item and its item skills. Minor fix to
# the skill tree front end." commit to the Items sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Minor fix to the skill tree front end: a handful of *_mod values on the
# existing artifacts (rows 2-5) were retuned.
# ---------------------------------------------------------------------

# Row 2 - Ancestral Finger Bone of The Magi Troth
$ws.Range("R2").Value = 0.25
$ws.Range("U2").Value = 0.1
$ws.Range("X2").Value = 0.4
$ws.Range("Z2").Value = 0.15

# Row 3 - Ancestral Witches Ice Bracelet
$ws.Range("U3").Value = 0.1
$ws.Range("W3").Value = 0.4
$ws.Range("Z3").Value = 0.15

# Row 4 - Emerald Laced Bow
$ws.Range("U4").Value = 0.1
$ws.Range("V4").Value = 0.4
$ws.Range("Y4").Value = 0.15

# Row 5 - Ancestral Soldiers Statue: id changed + a set of zero-valued
# stat columns were populated.
$ws.Range("A5").Value = 2406742
$ws.Range("O5").Value = 0
$ws.Range("Y5").Value = 0
$ws.Range("Z5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 0
$ws.Range("AK5").Value = 0
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0
$ws.Range("AS5").Value = 0
$ws.Range("AX5").Value = 0
$ws.Range("AY5").Value = 0
$ws.Range("BA5").Value = 0
$ws.Range("BB5").Value = 0
$ws.Range("BC5").Value = 0
$ws.Range("BD5").Value = 0
$ws.Range("BE5").Value = 0
$ws.Range("BF5").Value = 0
$ws.Range("BG5").Value = 0

# ---------------------------------------------------------------------
# Added a new ancestral item and its item skills: new artifact in row 6.
# ---------------------------------------------------------------------
$ws.Range("A6").Value = 3301314
$ws.Range("C6").Value = "Ancestral Fang of Delusional Thougts"
$ws.Range("D6").Value = "artifact"
$ws.Range("G6").Value = "A fang found in the memories of those deluded by the past. It contains the power of an ancient and yet powerful vampire: Tristie, she ruled the night and brought terror to those who crept through the shadows"
$ws.Range("Q6").Value = 0.5
$ws.Range("R6").Value = 0.25
$ws.Range("S6").Value = 0.3
$ws.Range("U6").Value = 0.4
$ws.Range("AC6").Value = 1
$ws.Range("AV6").Value = 0
$ws.Range("BM6").Value = 0
$ws.Range("BN6").Value = 0
$ws.Range("BO6").Value = 0
$ws.Range("BP6").Value = 0
$ws.Range("BQ6").Value = 0
$ws.Range("BT6").Value = "Twisted Blood Lust"

# ---------------------------------------------------------------------
# Minor column width tweaks (description column widened to fit the new
# longer text, and the item_skill_id column bumped to match).
# ---------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 244.3
$ws.Columns.Item(72).ColumnWidth = 21.5
